$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Country" used to be tracked on OrderProduct; it now belongs on Order, so
# it moves next to the other Order-level header fields: right after
# "Manufacturer" (col F) and before "Currency" (old col G). Insert a fresh
# column there, which pushes Currency..Update date one slot to the right
# and copies the left neighbour's (F1) formatting onto the new G1 header.
$ws.Range("G1").EntireColumn.Insert()

# Give the header its text - this also registers the new shared string.
$ws.Range("G1").Value = "Country"

# Size the new column the way it ended up after the edit.
$ws.Range("G1").ColumnWidth = 14

# Reflect the active cell/selection recorded after the edit.
$ws.Range("G13").Select()
